# eBayTestData.xlsx - "Final Commit With modifications"
#
# Adds three new columns (invalidUserName, invalidPassword, sortType) with
# sample values to the eBayTestData sheet, and makes that sheet the active /
# selected sheet (instead of DeviceCapabilities).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # eBayTestData
$ws2 = $wb.Worksheets.Item(2)   # DeviceCapabilities

# --- 1. Add the new header cells (row 1) and data cells (row 2) -----------
# The order in which new string values are first assigned controls the order
# they are appended to the shared-string table, so write them in the exact
# sequence that matches the target workbook (sortType ends up last).
$ws1.Range("N1").Value = "invalidUserName"
$ws1.Range("O1").Value = "invalidPassword"
$ws1.Range("N2").Value = "automation@gmail.com"
$ws1.Range("O2").Value = "Testing05"
$ws1.Range("P2").Value = "Highest Price"
$ws1.Range("P1").Value = "sortType"

# Match the bold header formatting used by the rest of row 1.
$ws1.Range("N1:P1").Font.Bold = $true

# --- 2. Size the three new columns similarly to the existing bestFit cols -
$ws1.Columns.Item(14).ColumnWidth = 15.6666666667
$ws1.Columns.Item(15).ColumnWidth = 14.6666666667
$ws1.Columns.Item(16).ColumnWidth = 12.5

# --- 3. Make eBayTestData the active sheet / tab, with P4 selected,
#        and make sure DeviceCapabilities is no longer the selected tab. --
$ws2.Range("B2").Select()
$ws1.Activate()
$ws1.Range("P4").Select()
